$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("soort vraag ...").
# This shifts the old column C ("totale score") -> D, and all question-score
# columns C..K -> D..L, copying formatting (style) from the left-hand column
# the same way Excel's own Insert does.
$ws.Columns("C:C").Insert()

# Set the new column's width (closest value this engine's quantised
# ColumnWidth can reach to the authored 15.6640625 stored width).
$ws.Columns("C:C").ColumnWidth = 14.75

# Header row labels
$ws.Range("C1").Value = "soort vraag (multiple = 0, één aantwoord mogeljk = 1)"
$ws.Range("D1").Value = "maximale score"

# New "soort vraag" column values (constant 1) for each question row
$ws.Range("C2:C9").Value = 1

# Row 10 (Q14) received new/updated score values, not merely the shifted
# old ones.
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 5
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = 1
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 1

# Restore the active selection recorded in the saved workbook
$null = $ws.Range("K13").Select()
